{"js": "// Replace each two-digit-division answer cell's text with its new value.\n// Source strings are unique in the document, so an exact-match search\n// (matchCase) for each one safely targets only the intended cell -\n// even though some NEW values repeat (e.g. \"64\\u00f72=32, 0\" is the\n// target of two different source cells), and even though the body also\n// contains a date heading paragraph that is left untouched.\nconst replacements = [\n  [\"88\u00f75=17, 3\", \"55\u00f77=7, 6\"],\n  [\"48\u00f75=9, 3\", \"23\u00f78=2, 7\"],\n  [\"54\u00f77=7, 5\", \"73\u00f74=18, 1\"],\n  [\"44\u00f76=7, 2\", \"10\u00f77=1, 3\"],\n  [\"36\u00f78=4, 4\", \"89\u00f74=22, 1\"],\n  [\"31\u00f76=5, 1\", \"21\u00f79=2, 3\"],\n  [\"58\u00f77=8, 2\", \"75\u00f75=15, 0\"],\n  [\"43\u00f79=4, 7\", \"77\u00f75=15, 2\"],\n  [\"82\u00f73=27, 1\", \"23\u00f74=5, 3\"],\n  [\"69\u00f73=23, 0\", \"73\u00f78=9, 1\"],\n  [\"39\u00f77=5, 4\", \"47\u00f78=5, 7\"],\n  [\"98\u00f76=16, 2\", \"27\u00f79=3, 0\"],\n  [\"68\u00f76=11, 2\", \"14\u00f73=4, 2\"],\n  [\"69\u00f79=7, 6\", \"60\u00f72=30, 0\"],\n  [\"37\u00f73=12, 1\", \"24\u00f79=2, 6\"],\n  [\"19\u00f72=9, 1\", \"12\u00f74=3, 0\"],\n  [\"85\u00f77=12, 1\", \"46\u00f74=11, 2\"],\n  [\"24\u00f78=3, 0\", \"94\u00f76=15, 4\"],\n  [\"42\u00f78=5, 2\", \"93\u00f73=31, 0\"],\n  [\"35\u00f72=17, 1\", \"17\u00f77=2, 3\"],\n  [\"10\u00f74=2, 2\", \"64\u00f72=32, 0\"],\n  [\"99\u00f79=11, 0\", \"64\u00f72=32, 0\"],\n  [\"26\u00f74=6, 2\", \"76\u00f73=25, 1\"],\n  [\"22\u00f78=2, 6\", \"51\u00f72=25, 1\"],\n  [\"76\u00f78=9, 4\", \"14\u00f73=4, 2\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Replace each two-digit-division answer cell's text with its new value.\n# Source strings are unique in the document (verified against the XML),\n# so a literal Find/Replace (MatchCase on, MatchWildcards off) for each\n# one safely targets only the intended cell, even though a couple of the\n# NEW values repeat (e.g. \"64\u00f72=32, 0\" lands in two different cells)\n# and the body also contains an untouched date-heading paragraph.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    (\"88\u00f75=17, 3\", \"55\u00f77=7, 6\"),\n    (\"48\u00f75=9, 3\", \"23\u00f78=2, 7\"),\n    (\"54\u00f77=7, 5\", \"73\u00f74=18, 1\"),\n    (\"44\u00f76=7, 2\", \"10\u00f77=1, 3\"),\n    (\"36\u00f78=4, 4\", \"89\u00f74=22, 1\"),\n    (\"31\u00f76=5, 1\", \"21\u00f79=2, 3\"),\n    (\"58\u00f77=8, 2\", \"75\u00f75=15, 0\"),\n    (\"43\u00f79=4, 7\", \"77\u00f75=15, 2\"),\n    (\"82\u00f73=27, 1\", \"23\u00f74=5, 3\"),\n    (\"69\u00f73=23, 0\", \"73\u00f78=9, 1\"),\n    (\"39\u00f77=5, 4\", \"47\u00f78=5, 7\"),\n    (\"98\u00f76=16, 2\", \"27\u00f79=3, 0\"),\n    (\"68\u00f76=11, 2\", \"14\u00f73=4, 2\"),\n    (\"69\u00f79=7, 6\", \"60\u00f72=30, 0\"),\n    (\"37\u00f73=12, 1\", \"24\u00f79=2, 6\"),\n    (\"19\u00f72=9, 1\", \"12\u00f74=3, 0\"),\n    (\"85\u00f77=12, 1\", \"46\u00f74=11, 2\"),\n    (\"24\u00f78=3, 0\", \"94\u00f76=15, 4\"),\n    (\"42\u00f78=5, 2\", \"93\u00f73=31, 0\"),\n    (\"35\u00f72=17, 1\", \"17\u00f77=2, 3\"),\n    (\"10\u00f74=2, 2\", \"64\u00f72=32, 0\"),\n    (\"99\u00f79=11, 0\", \"64\u00f72=32, 0\"),\n    (\"26\u00f74=6, 2\", \"76\u00f73=25, 1\"),\n    (\"22\u00f78=2, 6\", \"51\u00f72=25, 1\"),\n    (\"76\u00f78=9, 4\", \"14\u00f73=4, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
